$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (the "Förändrad" date) for rows 2 through 13 from 45212 to 45221
for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45212) {
        $cell.Value = 45221
    }
}
